# Updated cryptos list values (price + 1h volume change) per the source diff.
# Column D ("Price") holds plain numeric-looking strings (e.g. "0.998", "208.67") that
# must stay as TEXT, same as the original inline-string cells -- a bare assignment would
# let Excel coerce them into numbers (and drop trailing zeros). We force those with a
# leading apostrophe, then reset the cell style so no stray quote-prefix style lingers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "80.246.21" }
    @{ Cell = "E2"; Value = "  +5.05%  " }
    @{ Cell = "D3"; Value = "3.187.64" }
    @{ Cell = "E3"; Value = "  +3.60%  " }
    @{ Cell = "D4"; Value = "0.998" }
    @{ Cell = "E4"; Value = "  -0.17%  " }
    @{ Cell = "D5"; Value = "208.67" }
    @{ Cell = "E5"; Value = "  +5.26%  " }
    @{ Cell = "D6"; Value = "626.88" }
    @{ Cell = "E6"; Value = "  +1.07%  " }
    @{ Cell = "D7"; Value = "0.275" }
    @{ Cell = "E7"; Value = "  +28.31%  " }
    @{ Cell = "D8"; Value = "0.998" }
    @{ Cell = "E8"; Value = "  -0.14%  " }
    @{ Cell = "D9"; Value = "0.588" }
    @{ Cell = "E9"; Value = "  +6.47%  " }
    @{ Cell = "D10"; Value = "3.180.56" }
    @{ Cell = "E10"; Value = "  +3.39%  " }
    @{ Cell = "D11"; Value = "0.590" }
    @{ Cell = "E11"; Value = "  +30.94%  " }
    @{ Cell = "D12"; Value = "0.0000258" }
    @{ Cell = "E12"; Value = "  +28.98%  " }
    @{ Cell = "E13"; Value = "  +1.73%  " }
    @{ Cell = "D14"; Value = "3.766.13" }
    @{ Cell = "E14"; Value = "  +3.34%  " }
    @{ Cell = "D15"; Value = "5.26" }
    @{ Cell = "E15"; Value = "  +0.14%  " }
    @{ Cell = "D16"; Value = "31.86" }
    @{ Cell = "E16"; Value = "  +8.69%  " }
    @{ Cell = "D17"; Value = "79.775.17" }
    @{ Cell = "E17"; Value = "  +4.61%  " }
    @{ Cell = "D18"; Value = "3.171.81" }
    @{ Cell = "E18"; Value = "  +3.29%  " }
    @{ Cell = "D19"; Value = "14.26" }
    @{ Cell = "E19"; Value = "  +5.85%  " }
    @{ Cell = "D20"; Value = "3.00" }
    @{ Cell = "E20"; Value = "  +13.35%  " }
    @{ Cell = "D21"; Value = "9.14" }
    @{ Cell = "E21"; Value = "  +1.17%  " }
    @{ Cell = "D22"; Value = "436.43" }
    @{ Cell = "E22"; Value = "  +12.69%  " }
    @{ Cell = "D23"; Value = "5.16" }
    @{ Cell = "E23"; Value = "  +14.04%  " }
    @{ Cell = "D25"; Value = "3.335.18" }
    @{ Cell = "E25"; Value = "  +3.21%  " }
    @{ Cell = "D26"; Value = "75.98" }
    @{ Cell = "E26"; Value = "  +4.83%  " }
    @{ Cell = "E27"; Value = "  +2.31%  " }
    @{ Cell = "D28"; Value = "10.89" }
    @{ Cell = "E28"; Value = "  +6.68%  " }
    @{ Cell = "D29"; Value = "1.01" }
    @{ Cell = "E29"; Value = "  +0.62%  " }
    @{ Cell = "E30"; Value = "  +11.60%  " }
    @{ Cell = "E31"; Value = "  +0.36%  " }
    @{ Cell = "D32"; Value = "8.96" }
    @{ Cell = "E32"; Value = "  +7.60%  " }
    @{ Cell = "D33"; Value = "556.92" }
    @{ Cell = "E33"; Value = "  +11.00%  " }
    @{ Cell = "D34"; Value = "1.47" }
    @{ Cell = "E34"; Value = "  +2.42%  " }
    @{ Cell = "E35"; Value = "  +14.21%  " }
    @{ Cell = "D36"; Value = "1.99" }
    @{ Cell = "E36"; Value = "  +3.24%  " }
    @{ Cell = "D37"; Value = "22.98" }
    @{ Cell = "E37"; Value = "  +10.10%  " }
    @{ Cell = "D38"; Value = "0.123" }
    @{ Cell = "E38"; Value = "  +20.25%  " }
    @{ Cell = "E39"; Value = "  -0.07%  " }
    @{ Cell = "D40"; Value = "0.406" }
    @{ Cell = "E40"; Value = "  +7.59%  " }
    @{ Cell = "D41"; Value = "20.78" }
    @{ Cell = "E41"; Value = "  +3.52%  " }
    @{ Cell = "D42"; Value = "163.86" }
    @{ Cell = "E42"; Value = "  +0.16%  " }
    @{ Cell = "B43"; Value = "RenderToken" }
    @{ Cell = "C43"; Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render" }
    @{ Cell = "D43"; Value = "5.63" }
    @{ Cell = "E43"; Value = "  +8.38%  " }
    @{ Cell = "B44"; Value = "USDe" }
    @{ Cell = "C44"; Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde" }
    @{ Cell = "D44"; Value = "1.00" }
    @{ Cell = "E44"; Value = "  +0.00%  " }
    @{ Cell = "D45"; Value = "190.59" }
    @{ Cell = "E45"; Value = "  -1.73%  " }
    @{ Cell = "D46"; Value = "1.81" }
    @{ Cell = "E46"; Value = "  +8.60%  " }
    @{ Cell = "D47"; Value = "2.70" }
    @{ Cell = "E47"; Value = "  +9.33%  " }
    @{ Cell = "D48"; Value = "0.786" }
    @{ Cell = "E48"; Value = "  -1.54%  " }
    @{ Cell = "D49"; Value = "1.30" }
    @{ Cell = "E49"; Value = "  +2.72%  " }
    @{ Cell = "D50"; Value = "42.77" }
    @{ Cell = "E50"; Value = "  +4.35%  " }
    @{ Cell = "D51"; Value = "4.23" }
    @{ Cell = "E51"; Value = "  +7.72%  " }
)

# Cells whose new text would otherwise be auto-converted to a Number by Excel.
$forceTextCells = @(
    "D4", "D5", "D6", "D7", "D8", "D9", "D11", "D12", "D15", "D16", "D19", "D20", "D21", "D22", "D23", "D26", "D28", "D29", "D32", "D33", "D34", "D36", "D37", "D38", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51"
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)
    if ($forceTextCells -contains $u.Cell) {
        $range.Value = "'" + $u.Value
        $range.Style = "Normal"
    } else {
        $range.Value = $u.Value
    }
}
